$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.057685732841492
$ws.Range("B1").Value = 1.357896327972412
$ws.Range("C1").Value = 1.302277207374573
$ws.Range("D1").Value = 1.526670813560486
$ws.Range("E1").Value = 1.316855907440186
